$wb = $excel.ActiveWorkbook

# Sheet1 (展览)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 222
$ws1.Range("F5").Value = 1520
$ws1.Range("F6").Value = 212
$ws1.Range("F7").Value = 613
$ws1.Range("F8").Value = 128
$ws1.Range("F9").Value = 591
$ws1.Range("F10").Value = 43
$ws1.Range("F11").Value = 96
$ws1.Range("F13").Value = 157

# Sheet2 (演出)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 5
$ws2.Range("F11").Value = 11

# Sheet3 (本地生活)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1933

# Sheet4 (全部类型)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1933
$ws4.Range("F11").Value = 222
$ws4.Range("F12").Value = 5
$ws4.Range("F15").Value = 1520
$ws4.Range("F17").Value = 212
$ws4.Range("F19").Value = 11
$ws4.Range("F20").Value = 613
$ws4.Range("F22").Value = 128
$ws4.Range("F23").Value = 591
$ws4.Range("F24").Value = 43
$ws4.Range("F26").Value = 96
$ws4.Range("F30").Value = 157
